{"js": "// Fix Kieran Hughes' name (typo \"Kieran\" -> \"Keiran\") and student number\n// (\"91361\" -> \"913861\") in the contributions report.\nconst body = context.document.body;\n\nconst results = body.search(\"Kieran Hughes \u2013 91361\", { matchCase: true });\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find 'Kieran Hughes \u2013 91361' in the document.\");\n}\n\nresults.items[0].insertText(\"Keiran Hughes \u2013 913861\", \"Replace\");\nawait context.sync();\n", "ps1": "# Fix Kieran Hughes' name (typo \"Kieran\" -> \"Keiran\") and student number\n# (\"91361\" -> \"913861\") in the contributions report.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Kieran Hughes \u2013 91361\"\n$find.Forward = $true\n$find.Wrap = 1  # wdFindContinue\n\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Could not find 'Kieran Hughes \u2013 91361' in the document.\"\n}\n\n$find.Parent.Text = \"Keiran Hughes \u2013 913861\"\n"}
